# Auto-generated: applies 2024-01-16 daily crime data update across all sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range('K2').Value = 265
$ws.Range('K3').Value = 245
$ws.Range('K4').Value = 48
$ws.Range('K5').Value = 14
$ws.Range('K6').Value = 342
$ws.Range('K7').Value = 914

$ws = $wb.Worksheets.Item('Logan Square')
$ws.Range('K6').Value = 5
$ws.Range('K7').Value = 10

$ws = $wb.Worksheets.Item('Austin')
$ws.Range('K6').Value = 19
$ws.Range('K7').Value = 58

$ws = $wb.Worksheets.Item('South Chicago')
$ws.Range('K5').Value = 1
$ws.Range('K7').Value = 14

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range('K2').Value = 17
$ws.Range('K7').Value = 44

$ws = $wb.Worksheets.Item('West Pullman')
$ws.Range('K6').Value = 8
$ws.Range('K7').Value = 22

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range('K2').Value = 5
$ws.Range('K3').Value = 10
$ws.Range('K7').Value = 22

$ws = $wb.Worksheets.Item('New City')
$ws.Range('K2').Value = 7
$ws.Range('K3').Value = 6
$ws.Range('K5').Value = 1
$ws.Range('K7').Value = 24

$ws = $wb.Worksheets.Item('Woodlawn')
$ws.Range('K2').Value = 8
$ws.Range('K3').Value = 7
$ws.Range('K7').Value = 21

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range('K7').Value = 29
$ws.Range('K8').Value = 58
$ws.Range('K16').Value = 3
$ws.Range('K23').Value = 7
$ws.Range('K28').Value = 43
$ws.Range('K30').Value = 12
$ws.Range('K32').Value = 44
$ws.Range('K35').Value = 13
$ws.Range('K36').Value = 22
$ws.Range('K40').Value = 13
$ws.Range('K41').Value = 30
$ws.Range('K46').Value = 9
$ws.Range('K47').Value = 7
$ws.Range('K48').Value = 10
$ws.Range('K51').Value = 26
$ws.Range('K52').Value = 10
$ws.Range('K53').Value = 16
$ws.Range('K54').Value = 9
$ws.Range('K58').Value = 7
$ws.Range('K61').Value = 5
$ws.Range('K62').Value = 6
$ws.Range('K63').Value = 24
$ws.Range('K64').Value = 6
$ws.Range('K65').Value = 34
$ws.Range('K71').Value = 10
$ws.Range('K75').Value = 7
$ws.Range('K76').Value = 12
$ws.Range('K81').Value = 14
$ws.Range('K83').Value = 44
$ws.Range('K93').Value = 22
$ws.Range('K96').Value = 3
$ws.Range('K97').Value = 21
$ws.Range('K99').Value = 914

$ws = $wb.Worksheets.Item('Gage Park')
$ws.Range('K2').Value = 6
$ws.Range('K4').Value = 1
$ws.Range('K6').Value = 5
$ws.Range('K7').Value = 12

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range('K4').Value = 1
$ws.Range('K7').Value = 34

$ws = $wb.Worksheets.Item('Lincoln Park')
$ws.Range('K6').Value = 9
$ws.Range('K7').Value = 10

$ws = $wb.Worksheets.Item('Loop')
$ws.Range('K3').Value = 5
$ws.Range('K7').Value = 16

$ws = $wb.Worksheets.Item('Englewood')
$ws.Range('K2').Value = 14
$ws.Range('K3').Value = 12
$ws.Range('K6').Value = 16
$ws.Range('K7').Value = 43

$ws = $wb.Worksheets.Item('Lake View')
$ws.Range('K4').Value = 2
$ws.Range('K7').Value = 7

$ws = $wb.Worksheets.Item('Hermosa')
$ws.Range('K6').Value = 8
$ws.Range('K7').Value = 13

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range('K3').Value = 10
$ws.Range('K6').Value = 13
$ws.Range('K7').Value = 30

$ws = $wb.Worksheets.Item('Rogers Park')
$ws.Range('K5').Value = 6
$ws.Range('K6').Value = 12

$ws = $wb.Worksheets.Item('Lower West Side')
$ws.Range('K3').Value = 4
$ws.Range('K7').Value = 9

$ws = $wb.Worksheets.Item('Douglas')
$ws.Range('K5').Value = 4
$ws.Range('K6').Value = 7

$ws = $wb.Worksheets.Item('Near South Side')
$ws.Range('K3').Value = 3
$ws.Range('K7').Value = 6

$ws = $wb.Worksheets.Item('Grand Boulevard')
$ws.Range('K2').Value = 4
$ws.Range('K7').Value = 13

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Range('K2').Value = 10
$ws.Range('K3').Value = 9
$ws.Range('K7').Value = 29

$ws = $wb.Worksheets.Item('Kenwood')
$ws.Range('K3').Value = 3
$ws.Range('K7').Value = 9

$ws = $wb.Worksheets.Item('Wicker Park')
$ws.Range('K5').Value = 3
$ws.Range('K6').Value = 3

$ws = $wb.Worksheets.Item('North Center')
$ws.Range('J5').Value = 3
$ws.Range('J6').Value = 6

$ws = $wb.Worksheets.Item('Portage Park')
$ws.Range('K2').Value = 4
$ws.Range('K6').Value = 10

$ws = $wb.Worksheets.Item('Montclare')
$ws.Range('J6').Value = 4
$ws.Range('J7').Value = 7

$ws = $wb.Worksheets.Item('South Shore')
$ws.Range('K2').Value = 17
$ws.Range('K4').Value = 5
$ws.Range('K6').Value = 9
$ws.Range('K7').Value = 44

$ws = $wb.Worksheets.Item('Riverdale')
$ws.Range('K2').Value = 3
$ws.Range('K7').Value = 7

$ws = $wb.Worksheets.Item('Little Village')
$ws.Range('K2').Value = 6
$ws.Range('K3').Value = 6
$ws.Range('K7').Value = 26

$ws = $wb.Worksheets.Item('Bucktown')
$ws.Range('I4').Value = 2
$ws.Range('I5').Value = 3
